# "Switching to VLJ # in the spreadsheet"
#
# The "CSS Id" column is renamed to "VLJ #" and the example/sample rows are
# updated to use simple numeric VLJ numbers instead of the old CSS ids, plus
# the second judge's name is swapped out. Finally, a new (blank) row is
# appended at the bottom of the table, matching the formatting of the row
# above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "CSS Id" header to "VLJ #"
$ws.Range("C2").Value() = "VLJ #"

# Replace the sample CSS ids (BVAJONESB) with the new VLJ # values
$ws.Range("C3:C8").Value() = "123"

# The second example judge, "Roth, Lauren", becomes "Huels, Stuart"
$ws.Range("B8").Value() = "Huels, Stuart"
$ws.Range("B9").Value() = "Huels, Stuart"

# ...with her own distinct VLJ # on the last row
$ws.Range("C9").Value() = "456"

# Append a new blank row at the bottom of the table (row 10), copying the
# formatting (borders/fill/height) from the last existing row so it matches
# the rest of the table.
$ws.Range("A9:I9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$excel.CutCopyMode() = 0
$ws.Rows.Item(10).RowHeight() = 17
